# Sync attendance_reports: normalize "Recorded By" (column G) values so that
# "System" is listed first among the recorders for specific known value combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact literal replacements to apply to column G ("Recorded By") text values.
# Using an ordered map of old-value -> new-value, applied wherever the cell's
# text matches exactly.
$replacements = @{
    "system, backup@backdoor.com, System" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

# Determine the last used row from the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count
if ($ws.UsedRange.Row -gt 1) {
    $lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
